$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alldata")

# Apply an AutoFilter over the used range (A1:O245), filtering column G
# ("Owner", 7th column) down to just the rows where the value is "user".
# Passing an array for Criteria1 with Operator = xlFilterValues (7) makes
# Excel emit a discrete-value <filters><filter val="..."/></filters> list
# (rather than a <customFilters> comparison) and also hides the rows that
# don't match ("nico" rows), matching the recorded edit.
$range = $ws.Range("A1:O245")
$range.AutoFilter(7, @("user"), 7)

# Excel also records the active AutoFilter range as a hidden, sheet-scoped
# defined name (_xlnm._FilterDatabase).
$fdbName = $ws.Names.Add("_xlnm._FilterDatabase", "=alldata!`$A`$1:`$O`$245")
$fdbName.Visible = $false

# Move the view so row 158 is the top-most visible row, and select O2:O244
# with O2 as the active cell (matches the saved view state).
$ws.Application.ActiveWindow.ScrollRow = 158
$ws.Range("O2:O244").Select()

$wb.Save()
